$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iVals = @(8,9,8,8,9,7,6,6,9,8,4,8,9,8,7,6,6,8,7,7,9,6,7,7,5,6,9,7,8,6,6,9,7,9,7,7,8,8,9,8,7,8,7,8,7,7,8,10,6,10,7,6,7,6,9,9,7,8,7,7,7,7,7,6,6,7,7,7,7,5)
$jVals = @(8,9,8,8,10,8,7,7,9,9,5,9,9,8,7,8,7,8,8,8,9,7,7,8,6,7,9,7,8,6,6,9,8,9,8,7,9,9,9,8,7,8,7,8,8,7,8,10,7,11,7,7,7,6,9,9,8,8,7,7,7,7,7,7,7,7,8,8,7,5)

for ($r = 2; $r -le 71; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
